$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Urun_Ozellik_Bilgileri")

$ws.Range("B2").Value = "Apple Uyumlu"
$ws.Range("B4").Value = "Casper Uyumlu"
$ws.Range("B5").Value = "General Mobile Uyumlu"
$ws.Range("B6").Value = "Huawei Uyumlu"
$ws.Range("B7").Value = "JBL Uyumlu"
$ws.Range("B8").Value = "Lenovo Uyumlu"
$ws.Range("B9").Value = "Oppo Uyumlu"
$ws.Range("B10").Value = "POCO Uyumlu"
$ws.Range("B11").Value = "Reeder Uyumlu"
$ws.Range("B12").Value = "Samsung Uyumlu"
$ws.Range("B13").Value = "Sennheiser Uyumlu"
$ws.Range("B14").Value = "Sony Uyumlu"
$ws.Range("B15").Value = "Xiaomi Uyumlu"
